# "Digimon Ultimate Tier test"
# Adds three new skill rows (Paildramon's NewYearsEve / BackwaterCamp /
# DeathParadeBlaster set) at the bottom of the skills table, and updates
# the sheet's scroll position / selection to the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67 (A67:E67) already carries the "Bom" (good/green) cell style used
# by every data row in the table. Copy its formatting down into the three
# new rows (68:70) before stamping in the new values, so the new cells
# pick up the same style instead of Excel's default.
$ws.Range("A67:E67").Copy()
$ws.Range("A68:E70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 68 - id 66
$ws.Cells.Item(68, 1).Value = 66
$ws.Cells.Item(68, 2).Value = "NewYearsEve"
$ws.Cells.Item(68, 3).Value = "DamageSkill"
$ws.Cells.Item(68, 4).Value = 100
$ws.Cells.Item(68, 5).Value = 15

# Row 69 - id 67
$ws.Cells.Item(69, 1).Value = 67
$ws.Cells.Item(69, 2).Value = "BackwaterCamp"
$ws.Cells.Item(69, 3).Value = "PassiveSkill(Paildramon)"
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0

# Row 70 - id 68
$ws.Cells.Item(70, 1).Value = 68
$ws.Cells.Item(70, 2).Value = "DeathParadeBlaster"
$ws.Cells.Item(70, 3).Value = "DamageSkill(Paildramon)"
$ws.Cells.Item(70, 4).Value = 100
$ws.Cells.Item(70, 5).Value = 15

# Match the author's final view state: scrolled down to row 54, with C68
# selected.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("C68").Select()
